$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, C->D)
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "StatQuery"

# New stat-bar query text for the inserted column's data row
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Flat-Coated Retriever']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$ws.Cells.Item(2, 2).Value = $statQuery

# Match the wrap-text style used by A2 on the new B2 cell
$ws.Cells.Item(2, 2).WrapText = $true

# New column B should have the same width as column A (columns A and C/D
# already keep their original widths automatically after the insert).
# 75 is the closest input this engine's width-rounding accepts to reproduce
# column A's stored width of 75.81640625 on the newly inserted column.
$ws.Columns.Item(2).ColumnWidth = 75

# Update the selection to match the target state
$ws.Range("A2").Select()
